$d = $word.ActiveDocument

$replacements = @(
    @("92×43=", "30×58="),
    @("53×91=", "17×27="),
    @("75×20=", "11×47="),
    @("82×39=", "26×98="),
    @("32×64=", "57×81="),
    @("71×85=", "54×75="),
    @("93×43=", "74×92="),
    @("71×56=", "66×98="),
    @("53×31=", "98×29="),
    @("20×82=", "25×81="),
    @("45×19=", "34×75="),
    @("31×97=", "35×86="),
    @("47×75=", "95×76="),
    @("66×12=", "45×97="),
    @("93×24=", "98×98="),
    @("65×62=", "40×69="),
    @("63×89=", "99×21="),
    @("51×69=", "16×39="),
    @("54×79=", "58×25="),
    @("60×95=", "54×47="),
    @("65×66=", "29×55="),
    @("57×47=", "54×58="),
    @("86×56=", "40×36="),
    @("54×32=", "93×14="),
    @("66×19=", "53×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
